# Applies the reviewer-response edits described by the target diff:
#  1) "NHANES data." -> "NHANES dataset." (Response to Reviewer 1, comment 1)
#  2) "changed most of the results" -> "changed the results" (same paragraph)
#  3) Expands the "we adjusted for as many potential confounders..." sentence
#     with additional explanatory text about unmeasured confounders / the new
#     supplementary table with basic-adjustment models (Response, comment 2).
#  4) Removes the stale run split / lastRenderedPageBreak around "cancer
#     clinics." (now that new content pushed the real page break earlier in
#     the document, Word recomputes it) by normalizing the run through it.

$d = $word.ActiveDocument

# 1) "NHANES data" -> "NHANES dataset"
$d.Content.Find.Execute(
    "quite small in the NHANES data. Nevertheless",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "quite small in the NHANES dataset. Nevertheless", 2)

# 2) "changed most of the results" -> "changed the results"
$d.Content.Find.Execute(
    "Please note that this changed most of the results and, thus,",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Please note that this changed the results and, thus,", 2)

# 3) Expand the confounders sentence with new explanatory text
$d.Content.Find.Execute(
    "that we adjusted for as many potential confounders as we could. This is clearly indicated in the methods and results sections.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "that we adjusted for as many potential confounders as we could and there was a great deal of thought that went into the selection of the covariates to adjust for. This is clearly indicated in the methods and results sections. We cannot adjust for unmeasured confounders since they are unobserved. We have adjusted our analysis so that results for two other models (one that is unadjusted for any other confounders and a second that adjusted for a basic set of confounders" + [char]0x2014 + "namely age, sex, and race" + [char]0x2014 + "are now included in a separate supplementary table" + [char]0x2014 + "Supplementary Table SX). We also expanded the discussion of these results",
    2)

# 4) Normalize the run split around "cancer clinics." so the stale cached
#    lastRenderedPageBreak is dropped (a fresh one now lands earlier, at the
#    text inserted in step 3 above).
$d.Content.Find.Execute(
    "efforts to address this in cancer clinics. However",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "efforts to address this in cancer clinics. However", 2)
